# Mise a jour de l'application
#
# A new training/attendance day (2025-11-20) is appended as a new column
# right after the current last column (CC). The header row gets the new
# date serial, and every player row gets their attendance mark for that
# day copied from the previous day's column (CC) into the new column
# (CD) - matching the source diff exactly, mark-for-mark.
#
# Row 12's data does not extend as far as column CC (it stops at AX), so
# no cell is added there. Row 21's CC cell is blank (player inactive), so
# its new CD cell stays blank too, keeping only the formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new CD value ($null means "leave blank", only copy formatting)
$values = @{
    1  = 45981
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "B"
    6  = "B"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = $null
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($r in 1..29) {
    if ($r -eq 12) { continue }   # this row has no data out to column CC

    $dst = $ws.Cells.Item($r, 82)   # column CD
    $val = $values[$r]

    if ($null -ne $val) {
        $dst.Value2 = $val
    }

    # Copy the formatting (number format / alignment / style) from the
    # previous day's column (CC) into the new column, without disturbing
    # the value we just set (set value first, then bring over the format).
    $ws.Cells.Item($r, 81).Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats
}

$excel.CutCopyMode = 0

# Refresh the sheet's frozen-pane / active-cell view state to match the
# newly widened data range.
$ws.Range("A1").Select()
$pane = $ws.Panes.Item(1)
$pane.SplitColumn = 1
$ws.Application.ActiveWindow.FreezePanes = $true

$ws.Range("BZ1").Select()
$ws.Range("CF25").Select()
